# Applies the new "FPHR" functor-comment rows (52-54) to Sheet1,
# matching the commit "mapping of functors (some more comments)".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 52: Bloomingdale's example ------------------------------------
$ws.Range("B52").Value = "FPHR"
$ws.Range("E52").Value = "Mezi tyto řetězce patřily Bloomingdale's, vlastněný Campeau Corp. z Toro"
$ws.Range("E52").VerticalAlignment = -4160
$ws.Range("E52").HorizontalAlignment = -4131
$ws.Range("E52").WrapText = $true
$ws.Range("C52").Value = "slovo ""vlastněný"" by asi nemělo mít FPHR, ale třeba RSTR??"

# --- Row 53: Health Care Property Investors example --------------------
$ws.Range("B53").Value = "FPHR"
# E53 must stay on the default (unstyled) format, same as the sibling
# empty cells above it, so pull the format from one of those instead of
# letting a brand-new cell pick up column E's wrap-text default.
[void]$ws.Range("E21").Copy()
[void]$ws.Range("E53").PasteSpecial(-4122)
$ws.Range("E53").Value = "Health Care Property Investors Inc., nabízející 2250000 kmenových akcií"
$ws.Range("C53").Value = "slovo ""nabízející"" by asi nemělo mít FPHR, ale třeba RSTR??"

# --- Row 54: PS of New Hampshire example --------------------------------
$ws.Range("B54").Value = "FPHR"
# E54 needs the wrap-only style already used by E49.
[void]$ws.Range("E49").Copy()
[void]$ws.Range("E54").PasteSpecial(-4122)
$ws.Range("E54").Value = "… její předpovědi týkající se společnosti PS of New Hampshire - například růst poptávky po elektřině či zvýšená provozní efektivita - by se nevyplnily."
$ws.Range("C54").Value = "PAR ""například růst poptávky po elektřině či zvýšená provozní efektivita"" … má spíš viset na ""předpovědi"", možná jako RSTR"
$ws.Rows.Item(54).RowHeight = 45

# --- Selection / scroll position ---------------------------------------
[void]$ws.Range("B54").Select()

# --- Cosmetic: built-in style names were re-localised to English when
# the workbook was last saved by a newer Excel build (best effort; no-op
# if the host does not expose style renaming). ---------------------------
$wb.Styles.Item("Kontrolní buňka").Name = "Check Cell"
$wb.Styles.Item("Normální").Name = "Normal"

Write-Host "edit.ps1 applied"
